$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being treated as text, matching the
# original inline-string cell type, even for values that look numeric
# (e.g. "1.005", "0.00001309").
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{Row=2;  D="24.466.04"; E="  -1.47%  "},
    @{Row=3;  D="1.655.94";  E="  -3.03%  "},
    @{Row=4;  D="1.005";     E="  +0.32%  "},
    @{Row=5;  D="312.75";    E="  -0.87%  "},
    @{Row=6;  D="1.006";     E="  +0.44%  "},
    @{Row=7;  D="0.3925";    E="  -2.05%  "},
    @{Row=8;  D="0.3909";    E="  -3.31%  "},
    @{Row=9;  D="1.005";     E="  +0.29%  "},
    @{Row=10; D="50.69";     E="  -5.66%  "},
    @{Row=11; D="1.386";     E="  -5.96%  "},
    @{Row=12; D="0.08560";   E="  -2.85%  "},
    @{Row=13; D="25.06";     E="  -4.63%  "},
    @{Row=14; D="7.237";     E="  -3.83%  "},
    @{Row=15; D="0.00001309";E="  -2.61%  "},
    @{Row=16; D="7.618";     E="  -4.94%  "},
    @{Row=17; D="1.666.51";  E="  -2.41%  "},
    @{Row=18; D="93.24";     E="  -2.51%  "},
    @{Row=19; D="0.06990";   E="  -2.55%  "},
    @{Row=20; D="21.30";     E="  +1.90%  "},
    @{Row=21; D="7.007";     E="  -4.02%  "},
    @{Row=22; D="1.006";     E="  +0.50%  "},
    @{Row=23; D="13.83";     E="  -4.46%  "},
    @{Row=24; D="24.454.05"; E="  -1.51%  "},
    @{Row=25; D="2.345";     E="  -0.44%  "},
    @{Row=26; D="2.774";     E="  -4.40%  "},
    @{Row=27; D="22.79";     E="  -1.50%  "},
    @{Row=28; D="5.789";     E="  -6.52%  "},
    @{Row=29; D="158.95";    E="  -1.66%  "},
    @{Row=30; D="144.79";    E="  +0.47%  "},
    @{Row=31; D="8.232";     E="  +0.09%  "},
    @{Row=32; D="2.529";     E="  +10.94%  "},
    @{Row=33; D="1.843.60";  E="  -4.16%  "},
    @{Row=34; D="0.08207";   E="  -5.27%  "},
    @{Row=35; D="0.03035";   E="  -5.13%  "},
    @{Row=36; D="1.012";     E="  -1.78%  "},
    @{Row=37; D="6.865";     E="  -5.93%  "},
    @{Row=38; D="0.2778";    E="  -2.67%  "},
    @{Row=39; D="0.09576";   E="  +1.12%  "},
    @{Row=40; D="1.502";     E="  +1.31%  "},
    @{Row=41; D="10.23";     E="  -4.51%  "},
    @{Row=42; D="0.7813";    E="  -7.05%  "},
    @{Row=43; D="13.34";     E="  -6.29%  "},
    @{Row=44; D="16.41";     E="  -6.65%  "},
    @{Row=45; D="2.561";     E="  -5.94%  "},
    @{Row=46; D="0.7034";    E="  -5.24%  "},
    @{Row=47; D="4.156";     E="  -1.48%  "},
    @{Row=48; D="1.005";     E="  +0.34%  "},
    @{Row=49; D="0.08570";   E="  +2.01%  "},
    @{Row=50; D="1.312";     E="  -5.21%  "},
    @{Row=51; D="137.45";    E="  -2.23%  "}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value2 = $u.D
    $ws.Cells.Item($u.Row, 5).Value2 = $u.E
}
